# Update "想去人数" (column F) values across the worksheets to reflect
# freshly scraped counts (gh-pages output regenerated at 456a3b4).

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F3").Value = 3903
$wsExhibit.Range("F4").Value = 2308
$wsExhibit.Range("F5").Value = 457
$wsExhibit.Range("F10").Value = 113
$wsExhibit.Range("F12").Value = 256
$wsExhibit.Range("F13").Value = 2586
$wsExhibit.Range("F14").Value = 179

# --- Sheet "演出" ---
$wsShow = $wb.Worksheets.Item("演出")
$wsShow.Range("F2").Value = 36

# --- Sheet "本地生活" --- (no changes)

# --- Sheet "全部类型" ---
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F3").Value = 3903
$wsAll.Range("F4").Value = 2308
$wsAll.Range("F5").Value = 457
$wsAll.Range("F8").Value = 36
$wsAll.Range("F11").Value = 113
$wsAll.Range("F15").Value = 256
$wsAll.Range("F16").Value = 2586
$wsAll.Range("F17").Value = 179
